$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.788.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.48%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.292.41"
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.79%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.11%  "

# Row 9
$ws.Range("E9").Value = "  +3.65%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.89%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.47%  "

# Row 12
$ws.Range("E12").Value = "  +1.83%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.60%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.50%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.651.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.298.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.55%  "

# Row 17
$ws.Range("E17").Value = "  +0.63%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.700.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.50%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.81%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0898"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.19%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.69%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.55%  "

# Row 25
$ws.Range("E25").Value = "  -0.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.64%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.67%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.63%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.38%  "

# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.79%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.45%  "

# Row 32
$ws.Range("E32").Value = "  -0.07%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.36%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.88%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.19%  "

# Row 36
$ws.Range("E36").Value = "  +2.81%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0689"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.28%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.101"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.82%  "

# Row 39
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "

# Row 40
$ws.Range("E40").Value = "  +1.61%  "

# Row 41
$ws.Range("E41").Value = "  -0.29%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.986.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.04%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0288"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.24%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.63%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.09%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.17%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.00%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.79%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.518.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.23%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.51%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.26%  "
